# Fruta / hortaliza, semanal
# Re-order the weekly price rows (2-8) by date, moving the D/M/N/O/P/S
# values of each row to their corresponding new row while leaving the
# descriptive columns (which are identical across rows) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) values for the columns that move between rows.
$cols = @("D", "M", "N", "O", "P", "S")
$rows = 2..8

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Mapping: new row -> source (old) row, derived from the target diff.
$mapping = @{
    2 = 5
    3 = 2
    4 = 3
    5 = 4
    6 = 8
    7 = 7
    8 = 6
}

foreach ($newRow in $rows) {
    $srcRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $orig[$srcRow][$c]
    }
}
